$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rule row for "R30" (18-21) had its greeting changed from
# "Good Evening" to "Good" in cell E10.
$ws.Range("E10").Value = "Good"
